# Auto-generated Excel COM-interop script to apply symbol-list / price update
# Commit: Updated symbol list on Tue Dec 20 08:16:28 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# For every touched cell, force Text number format first so that
# numeric-looking values (prices, hour counters) are written back
# as text strings (matching the inlineStr storage used by the source data)
# rather than being auto-coerced into numeric cell values.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '248.22'
$ws.Range('G2').NumberFormat = "@"
$ws.Range('G2').Value = '8'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '21.74'
$ws.Range('G3').NumberFormat = "@"
$ws.Range('G3').Value = '8'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.341'
$ws.Range('G4').NumberFormat = "@"
$ws.Range('G4').Value = '8'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05623'
$ws.Range('G5').NumberFormat = "@"
$ws.Range('G5').Value = '8'
$ws.Range('G6').NumberFormat = "@"
$ws.Range('G6').Value = '8'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.372'
$ws.Range('G7').NumberFormat = "@"
$ws.Range('G7').Value = '8'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8166'
$ws.Range('G8').NumberFormat = "@"
$ws.Range('G8').Value = '8'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9489'
$ws.Range('G9').NumberFormat = "@"
$ws.Range('G9').Value = '8'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1430'
$ws.Range('G10').NumberFormat = "@"
$ws.Range('G10').Value = '8'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07586'
$ws.Range('G11').NumberFormat = "@"
$ws.Range('G11').Value = '8'
$ws.Range('G12').NumberFormat = "@"
$ws.Range('G12').Value = '8'
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'ProBitToken'
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.1311'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '12ProBitTokenPROB'
$ws.Range('G13').NumberFormat = "@"
$ws.Range('G13').Value = '8'
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.03091'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('G14').NumberFormat = "@"
$ws.Range('G14').Value = '8'
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.09309'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('G15').NumberFormat = "@"
$ws.Range('G15').Value = '8'
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.561'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('G16').NumberFormat = "@"
$ws.Range('G16').Value = '8'
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.001594'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('G17').NumberFormat = "@"
$ws.Range('G17').Value = '8'
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.04699'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('G18').NumberFormat = "@"
$ws.Range('G18').Value = '8'
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'One'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0005781'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '18OneONE'
$ws.Range('G19').NumberFormat = "@"
$ws.Range('G19').Value = '8'
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'TigerCash'
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.006323'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '19TigerCashTCH'
$ws.Range('G20').NumberFormat = "@"
$ws.Range('G20').Value = '8'
$ws.Range('B21').NumberFormat = "@"
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').NumberFormat = "@"
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.005059'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('G21').NumberFormat = "@"
$ws.Range('G21').Value = '8'
$ws.Range('B22').NumberFormat = "@"
$ws.Range('B22').Value = 'BitKan'
$ws.Range('C22').NumberFormat = "@"
$ws.Range('C22').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.001034'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '21BitKanKAN'
$ws.Range('G22').NumberFormat = "@"
$ws.Range('G22').Value = '8'
$ws.Range('B23').NumberFormat = "@"
$ws.Range('B23').Value = 'NitroEx'
$ws.Range('C23').NumberFormat = "@"
$ws.Range('C23').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.0001501'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '22NitroExNTX'
$ws.Range('G23').NumberFormat = "@"
$ws.Range('G23').Value = '8'
$ws.Range('B24').NumberFormat = "@"
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').NumberFormat = "@"
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.769'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('G24').NumberFormat = "@"
$ws.Range('G24').Value = '8'
$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.141'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('G25').NumberFormat = "@"
$ws.Range('G25').Value = '8'
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.3302'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('G26').NumberFormat = "@"
$ws.Range('G26').Value = '8'
$ws.Range('G27').NumberFormat = "@"
$ws.Range('G27').Value = '8'
$ws.Range('G28').NumberFormat = "@"
$ws.Range('G28').Value = '8'
$ws.Range('G29').NumberFormat = "@"
$ws.Range('G29').Value = '8'
$ws.Range('G30').NumberFormat = "@"
$ws.Range('G30').Value = '8'
$ws.Range('G31').NumberFormat = "@"
$ws.Range('G31').Value = '8'
$ws.Range('G32').NumberFormat = "@"
$ws.Range('G32').Value = '8'
$ws.Range('G33').NumberFormat = "@"
$ws.Range('G33').Value = '8'
$ws.Range('G34').NumberFormat = "@"
$ws.Range('G34').Value = '8'
$ws.Range('G35').NumberFormat = "@"
$ws.Range('G35').Value = '8'
$ws.Range('G36').NumberFormat = "@"
$ws.Range('G36').Value = '8'
$ws.Range('G37').NumberFormat = "@"
$ws.Range('G37').Value = '8'
$ws.Range('G38').NumberFormat = "@"
$ws.Range('G38').Value = '8'
$ws.Range('G39').NumberFormat = "@"
$ws.Range('G39').Value = '8'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03959'
$ws.Range('G40').NumberFormat = "@"
$ws.Range('G40').Value = '8'
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1067'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('G41').NumberFormat = "@"
$ws.Range('G41').Value = '8'
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.003032'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('G42').NumberFormat = "@"
$ws.Range('G42').Value = '8'
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002925'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('G43').NumberFormat = "@"
$ws.Range('G43').Value = '8'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.008825'
$ws.Range('G44').NumberFormat = "@"
$ws.Range('G44').Value = '8'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005579'
$ws.Range('G45').NumberFormat = "@"
$ws.Range('G45').Value = '8'
$ws.Range('G46').NumberFormat = "@"
$ws.Range('G46').Value = '8'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '46ACDXExchangeACXTWorstin24h'
$ws.Range('G47').NumberFormat = "@"
$ws.Range('G47').Value = '8'
$ws.Range('G48').NumberFormat = "@"
$ws.Range('G48').Value = '8'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1786'
$ws.Range('G49').NumberFormat = "@"
$ws.Range('G49').Value = '8'
$ws.Range('G50').NumberFormat = "@"
$ws.Range('G50').Value = '8'
$ws.Range('G51').NumberFormat = "@"
$ws.Range('G51').Value = '8'
